$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift all existing data rows (2-119) down by one row (to 3-120),
# preserving values and existing per-cell formatting (e.g. the date style on column D).
$src = $ws.Range("A2:R119")
$dst = $ws.Range("A3:R120")
$dst.Value2 = $src.Value2

# The newly created row 120 needs the same date display format as the rest of column D.
$ws.Range("D120").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Write this week's new price data into row 2.
$ws.Range("D2").Value2 = 44956
$ws.Range("J2").Value2 = 430
$ws.Range("K2").Value2 = 23000
$ws.Range("L2").Value2 = 25000
$ws.Range("M2").Value2 = 24000
$ws.Range("P2").Value2 = 1846
